# Weekly update: insert a new price record as row 290, pushing the
# existing rows 290-388 down to 291-389 (dimension grows to A1:T389).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 290; everything below shifts down.
$ws.Rows.Item(290).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(290, 1).Value = 6
$ws.Cells.Item(290, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(290, 3).Value = 'Metropolitana'
$ws.Cells.Item(290, 4).Value = 45093
$ws.Cells.Item(290, 5).Value = 13
$ws.Cells.Item(290, 6).Value = 'Fruta'
$ws.Cells.Item(290, 7).Value = 100101
$ws.Cells.Item(290, 8).Value = 'Berries'
$ws.Cells.Item(290, 9).Value = 100101004
$ws.Cells.Item(290, 10).Value = 'Frambuesa'
$ws.Cells.Item(290, 11).Value = 'Sin especificar'
$ws.Cells.Item(290, 12).Value = 'Primera'
$ws.Cells.Item(290, 13).Value = 50
$ws.Cells.Item(290, 14).Value = 10000
$ws.Cells.Item(290, 15).Value = 10000
$ws.Cells.Item(290, 16).Value = 10000
$ws.Cells.Item(290, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(290, 18).Value = 'Región del Maule'
$ws.Cells.Item(290, 19).Value = 5000
$ws.Cells.Item(290, 20).Value = 2
